$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("urto inclinato")

# Dati raccolti per il piano inclinato (colonne C = "m2 a", D = "m2 sigmaa")

$ws.Range("C2").Value = 0.393
$ws.Range("D2").Value = 0.0013

$ws.Range("C3").Value = 0.394
$ws.Range("D3").Value = 0.0013

$ws.Range("C4").Value = 0.389
$ws.Range("D4").Value = 0.0011

$ws.Range("C5").Value = 0.388
$ws.Range("D5").Value = 0.0014

$ws.Range("C6").Value = 0.385
$ws.Range("D6").Value = 0.0012

$ws.Range("C7").Value = 0.385
$ws.Range("D7").Value = 0.001

# Righe 8 e 9: valori incollati da un'altra origine, formato ripristinato a
# generale e celle bloccate esplicitamente
$ws.Range("C8").Value = 0.385
$ws.Range("D8").Value = 0.0006
$ws.Range("C9").Value = 0.384
$ws.Range("D9").Value = 0.0008
$ws.Range("C8:D9").NumberFormat = "General"
$ws.Range("C8:D9").Locked = $true

$ws.Range("C10").Value = 0.399
$ws.Range("D10").Value = 0.0008

$ws.Range("C11").Value = 0.387
$ws.Range("D11").Value = 0.0009
$ws.Range("D11").NumberFormat = "0.0000"
$ws.Range("D11").Locked = $true

$ws.Range("C12").Value = 0.385
$ws.Range("D12").Value = 0.0008

# due righe vuote aggiunte in fondo alla tabella
$ws.Rows.Item(13).RowHeight = 14.25
$ws.Rows.Item(14).RowHeight = 14.25
